# Data-driven update of DAMSLTag (col I) and DialogAct (col J) for specific rows
$updates = @(
    @{Row=2; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=4; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=14; I='ba'; J='Appreciation'}
    @{Row=19; I='sv'; J='Statement-opinion'}
    @{Row=21; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=22; I='sd'; J='Statement-non-opinion'}
    @{Row=23; I='sv'; J='Statement-opinion'}
    @{Row=54; I='aa'; J='Agree/Accept'}
    @{Row=84; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=88; I='sd'; J='Statement-non-opinion'}
    @{Row=103; I='sd'; J='Statement-non-opinion'}
    @{Row=111; I='sd'; J='Statement-non-opinion'}
    @{Row=116; I='%'; J='Uninterpretable'}
    @{Row=117; I='%'; J='Uninterpretable'}
    @{Row=119; I='sd'; J='Statement-non-opinion'}
    @{Row=133; I='sv'; J='Statement-opinion'}
    @{Row=139; I='sv'; J='Statement-opinion'}
    @{Row=140; I='sv'; J='Statement-opinion'}
    @{Row=143; I='sd'; J='Statement-non-opinion'}
    @{Row=144; I='sd'; J='Statement-non-opinion'}
    @{Row=156; I='sv'; J='Statement-opinion'}
    @{Row=157; I='ba'; J='Appreciation'}
    @{Row=159; I='sd'; J='Statement-non-opinion'}
    @{Row=160; I='sv'; J='Statement-opinion'}
    @{Row=161; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=164; I='sd'; J='Statement-non-opinion'}
    @{Row=166; I='sd'; J='Statement-non-opinion'}
    @{Row=172; I='sd'; J='Statement-non-opinion'}
    @{Row=190; I='sd'; J='Statement-non-opinion'}
    @{Row=200; I='sd'; J='Statement-non-opinion'}
    @{Row=203; I='sv'; J='Statement-opinion'}
    @{Row=207; I='aa'; J='Agree/Accept'}
    @{Row=223; I='qy'; J='Yes-No-Question'}
    @{Row=227; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=257; I='ba'; J='Appreciation'}
    @{Row=271; I='sd'; J='Statement-non-opinion'}
    @{Row=283; I='aa'; J='Agree/Accept'}
    @{Row=284; I='aa'; J='Agree/Accept'}
    @{Row=288; I='%'; J='Uninterpretable'}
    @{Row=289; I='sv'; J='Statement-opinion'}
    @{Row=295; I='sd'; J='Statement-non-opinion'}
    @{Row=318; I='sd'; J='Statement-non-opinion'}
    @{Row=322; I='sv'; J='Statement-opinion'}
    @{Row=326; I='sv'; J='Statement-opinion'}
    @{Row=330; I='ba'; J='Appreciation'}
    @{Row=347; I='sd'; J='Statement-non-opinion'}
    @{Row=356; I='sv'; J='Statement-opinion'}
    @{Row=377; I='sd'; J='Statement-non-opinion'}
    @{Row=379; I='%'; J='Uninterpretable'}
    @{Row=383; I='sv'; J='Statement-opinion'}
    @{Row=395; I='aa'; J='Agree/Accept'}
    @{Row=399; I='sd'; J='Statement-non-opinion'}
    @{Row=403; I='ba'; J='Appreciation'}
    @{Row=411; I='sd'; J='Statement-non-opinion'}
    @{Row=438; I='sv'; J='Statement-opinion'}
    @{Row=445; I='sd'; J='Statement-non-opinion'}
    @{Row=454; I='sv'; J='Statement-opinion'}
    @{Row=456; I='sv'; J='Statement-opinion'}
)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}
